# edit.ps1 -- apply the two changes described by the diff:
#  1. "NOTAS IMPORTANTES:" paragraph: make bold+red and bump size 12pt -> 18pt
#  2. Last paragraph: merge the two runs that were split around the
#     "_GoBack" bookmark into a single run with the full sentence, keeping
#     the (now-collapsed) bookmark positioned right after the merged run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "NOTAS IMPORTANTES:" heading -> bold, red (FF0000), size 18pt (sz 36)
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*NOTAS IMPORTANTES*") {
        $p.Range.Font.Color = 255
        $p.Range.Font.Size = 18
    }
}

# ---------------------------------------------------------------------
# 2) Merge the split runs in the final paragraph around the _GoBack
#    bookmark, keeping the bookmark collapsed right after the merged run.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")

    # Find which paragraph holds the bookmark, and its end (excluding the
    # paragraph mark character).
    $targetPara = $null
    foreach ($p in $d.Paragraphs) {
        if (($bm.Start -ge $p.Range.Start) -and ($bm.Start -le $p.Range.End)) {
            $targetPara = $p
        }
    }

    if ($targetPara -ne $null) {
        $paraTextEnd = $targetPara.Range.End - 1

        # Range that holds the text after the bookmark up to (not
        # including) the paragraph mark -- this is the "second run".
        $afterRange = $d.Range($bm.End, $paraTextEnd)
        $afterText = $afterRange.Text

        if ($afterText -ne "") {
            # Remove that trailing text.
            $afterRange.Text = ""

            # Re-fetch the (still-collapsed) bookmark and insert the
            # removed text right at its position; Word keeps the
            # collapsed bookmark anchored *after* newly inserted text,
            # so the bookmark ends up right after the merged run again,
            # now sitting at the (new) end of the paragraph.
            $bm2 = $d.Bookmarks("_GoBack")
            $insertionPoint = $d.Range($bm2.Start, $bm2.Start)
            $insertionPoint.InsertAfter($afterText)
        }
    }
}
